# PPPQ Power BI Summit quiz — update the cover slide's edition label.
#
# Slide 1 ("Power Platform / Pub Quiz" cover), shape id=13 ("Rectangle 12")
# previously read "#MVP EDITION"; it becomes "#POWER BI SUMMIT" for the
# Power BI Summit edition of the quiz.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Rectangle 12")
$shape.TextFrame.TextRange.Runs(1).Text = "#POWER BI SUMMIT"
